# Update countries & provincias Spain
# Applies the Aug 9 2020 16:02 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp update -------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 9 de Agosto de 2020 a las 16:02"

# --- Per-country statistic refresh --------------------------------------
# Row 4: Estados Unidos
$ws.Range("B4").Value = 5152001
$ws.Range("C4").Value = 2278
$ws.Range("D4").Value = 2638713
$ws.Range("E4").Value = 2348195
$ws.Range("G4").Value = 23
$ws.Range("H4").Value = 165093

# Row 6: India
$ws.Range("B6").Value = 2167420
$ws.Range("C6").Value = 15400
$ws.Range("D6").Value = 1489005
$ws.Range("E6").Value = 634832
$ws.Range("G6").Value = 130
$ws.Range("H6").Value = 43583

# Row 15: Reino Unido
$ws.Range("B15").Value = 310825
$ws.Range("C15").Value = 1062

# Row 20: Argentina
$ws.Range("E20").Value = 129013
$ws.Range("G20").Value = 33
$ws.Range("H20").Value = 4556

# Row 22: Alemania
$ws.Range("B22").Value = 216911
$ws.Range("C22").Value = 15
$ws.Range("E22").Value = 10250

# Row 24: Irak
$ws.Range("B24").Value = 150115
$ws.Range("C24").Value = 2726
$ws.Range("D24").Value = 107775
$ws.Range("E24").Value = 36948
$ws.Range("G24").Value = 82
$ws.Range("H24").Value = 5392

# Row 62: Uzbekistan
$ws.Range("B62").Value = 30464
$ws.Range("C62").Value = 812
$ws.Range("D62").Value = 21813
$ws.Range("E62").Value = 8459
$ws.Range("G62").Value = 5
$ws.Range("H62").Value = 192

# Row 63: Serbia
$ws.Range("B63").Value = 28099
$ws.Range("C63").Value = 236
$ws.Range("E63").Value = 13411
$ws.Range("G63").Value = 9
$ws.Range("H63").Value = 641

# Row 86: Noruega
$ws.Range("B86").Value = 9615
$ws.Range("C86").Value = 16
$ws.Range("E86").Value = 502

# Row 128: Eslovenia
$ws.Range("B128").Value = 2249
$ws.Range("C128").Value = 2
$ws.Range("E128").Value = 195
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 127

# Row 170: Birmania
$ws.Range("B170").Value = 360
$ws.Range("C170").Value = 1
$ws.Range("D170").Value = 312

# Rows 180-184: small-territory block refresh
# Row 180: Islas Caimanes
$ws.Range("B180").Value = 214
$ws.Range("C180").Value = 26
$ws.Range("D180").Value = 53
$ws.Range("E180").Value = 158
$ws.Range("H180").Value = 3

# Row 181: Islas Turcas y Caicos
$ws.Range("B181").Value = 203
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 202
$ws.Range("E181").Value = 0
$ws.Range("H181").Value = 1

# Row 182: Gibraltar
$ws.Range("C182").Value = 27
$ws.Range("D182").Value = 39
$ws.Range("E182").Value = 156
$ws.Range("H182").Value = 2

# Row 183: San Martin (Parte Holandesa)
$ws.Range("B183").Value = 197
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 184
$ws.Range("E183").Value = 13
$ws.Range("H183").Value = 0

# Row 184: Papua Nueva Guinea
$ws.Range("B184").Value = 189
$ws.Range("C184").Value = 12
$ws.Range("D184").Value = 86
$ws.Range("E184").Value = 86
$ws.Range("H184").Value = 17
